$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "34.763.33"
Set-TextValue "E2" "  +0.96%  "
Set-TextValue "D3" "1.820.16"
Set-TextValue "E3" "  +1.31%  "
Set-TextValue "E4" "  -0.23%  "
Set-TextValue "D5" "228.75"
Set-TextValue "E5" "  +0.85%  "
Set-TextValue "D6" "0.578"
Set-TextValue "E6" "  +4.19%  "
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "D8" "34.89"
Set-TextValue "E8" "  +7.25%  "
Set-TextValue "D9" "0.301"
Set-TextValue "E9" "  +1.76%  "
Set-TextValue "D10" "0.0700"
Set-TextValue "E10" "  +1.02%  "
Set-TextValue "D11" "0.0954"
Set-TextValue "E11" "  +0.34%  "
Set-TextValue "D12" "2.081.23"
Set-TextValue "E12" "  +1.24%  "
Set-TextValue "D13" "11.52"
Set-TextValue "E13" "  +3.95%  "
Set-TextValue "D14" "1.828.51"
Set-TextValue "E14" "  +1.90%  "
Set-TextValue "E15" "  +2.48%  "
Set-TextValue "D16" "34.757.17"
Set-TextValue "E16" "  +1.07%  "
Set-TextValue "D17" "4.35"
Set-TextValue "E17" "  +2.89%  "
Set-TextValue "D18" "69.11"
Set-TextValue "D19" "248.41"
Set-TextValue "E19" "  +0.62%  "
Set-TextValue "D20" "0.0₃0805"
Set-TextValue "E20" "  +0.32%  "
Set-TextValue "D21" "11.60"
Set-TextValue "E21" "  +5.61%  "
Set-TextValue "E22" "  -0.17%  "
Set-TextValue "E23" "  +1.64%  "
Set-TextValue "D24" "171.81"
Set-TextValue "E24" "  +5.72%  "
Set-TextValue "E25" "  +1.81%  "
Set-TextValue "D26" "7.47"
Set-TextValue "E26" "  +3.83%  "
Set-TextValue "D27" "16.80"
Set-TextValue "E27" "  +2.34%  "
Set-TextValue "D28" "0.119"
Set-TextValue "E28" "  +2.80%  "
Set-TextValue "E29" "  -0.40%  "
Set-TextValue "D30" "3.99"
Set-TextValue "E30" "  +2.80%  "
Set-TextValue "E31" "  +2.27%  "
Set-TextValue "D32" "3.86"
Set-TextValue "E32" "  +2.11%  "
Set-TextValue "E33" "  +0.78%  "
Set-TextValue "E34" "  +1.44%  "
Set-TextValue "E35" "  +0.31%  "
Set-TextValue "D36" "1.417.36"
Set-TextValue "E36" "  -1.90%  "
Set-TextValue "E37" "  +2.10%  "
Set-TextValue "E38" "  +1.87%  "
Set-TextValue "E39" "  +0.95%  "
Set-TextValue "D40" "85.39"
Set-TextValue "E40" "  +1.78%  "
Set-TextValue "E41" "  +3.67%  "
Set-TextValue "D42" "0.959"
Set-TextValue "E42" "  +2.54%  "
Set-TextValue "E43" "  +0.11%  "
Set-TextValue "E44" "  +0.08%  "
Set-TextValue "E45" "  +3.35%  "
Set-TextValue "E46" "  -1.34%  "
Set-TextValue "D47" "6.11"
Set-TextValue "E47" "  +0.26%  "
Set-TextValue "D48" "1.982.65"
Set-TextValue "E48" "  +1.58%  "
Set-TextValue "D49" "105.59"
Set-TextValue "E49" "  -0.04%  "
Set-TextValue "E50" "  +0.18%  "
Set-TextValue "E51" "  -0.20%  "
